# Add User Stories for the new "Social Networks" domain inside the
# "Social and Urban Studies" cluster (mirrors the existing "Social Media"
# domain block, rows 82-101, copied down to rows 102-121).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Data: the 20 Machine Learning Task labels (col D) in the fixed order
# already used throughout the sheet, and the 20 new User Story texts
# (col E) for the "Social Networks" domain, in the same order.
# ---------------------------------------------------------------------
$tasks = @(
  'adversarial learning'
  'cnn'
  'conversational agent'
  'decision tree'
  'document classification'
  'entity extraction'
  'feature selection'
  'imbalanced dataset'
  'k-nearest neighbor'
  'keyword extraction'
  'multi-label classification'
  'neural network'
  'random forest'
  'semantic similarity'
  'sentiment analysis'
  'speech to text'
  'text categorization'
  'unsupervised clustering'
  'voice recognition'
  'word embedding'
)
$stories = @(
  'As a social network analyst, I want to apply adversarial learning techniques to detect and prevent fake accounts and bot activities, so that I can maintain the integrity and trustworthiness of the platform.'
  'As a social network content curator, I want to use CNNs to analyze video content for detecting trending topics and sentiment analysis among users, so that I can optimize content recommendations and advertising strategies.'
  'As a content curator for a social networking site, I want to deploy a conversational agent that suggests personalized content (articles, videos, posts) to users based on their interests and engagement history, so that I can enhance user engagement and retention.'
  'As a social network analyst, I want to build decision tree models to analyze user behavior patterns (e.g., posting frequency, interaction preferences), so that I can understand user engagement trends and optimize platform algorithms.'
  'As a community manager, I want to deploy document classification models to detect and categorize mentions of events or activities in user posts and comments, so that I can facilitate event promotion and engagement among users.'
  'As a marketing strategist, I want to apply entity extraction to identify and classify influential users and celebrities mentioned in social media conversations, so that I can engage with key influencers for promotional campaigns and collaborations.'
  'As a trend researcher, I want to utilize feature selection techniques to extract relevant features (e.g., hashtag usage, topic frequency) from social media posts to uncover emerging trends and topics of interest among users, so that I can provide insights for content creation and marketing strategies.'
  'As a sentiment analyst, I want to manage imbalanced sentiment datasets (e.g., disproportionately positive or negative comments) on social networks to develop accurate sentiment analysis models, so that I can understand nuanced public opinion and sentiment trends.'
  'As a social network analyst, I want to use KNN algorithms to identify clusters of users based on their communication patterns (e.g., messaging frequency, content similarity), so that I can understand community structures and interactions within the platform.'
  "As a marketing strategist, I want to apply keyword extraction techniques to analyze competitors' social media content and extract key themes and strategies they employ, so that I can benchmark and refine our own marketing efforts accordingly."
  "As a marketing strategist, I want to employ multi-label classification to profile users' interests across various categories (e.g., travel, food, music) from their social media activities, so that I can target them with relevant advertisements and promotions."
  'As a content moderator, I want to deploy neural networks for deep learning-based content filtering to automatically detect and flag inappropriate or sensitive content (e.g., nudity, violence) on social media platforms, so that I can maintain a safe and positive user experience.'
  'As a social media platform manager, I want to implement a random forest model to detect fake accounts based on multiple behavioral and profile characteristics (e.g., posting frequency, account creation details), so that I can enhance platform security and user trust.'
  'As a search engine developer, I want to leverage semantic similarity techniques to expand user queries by identifying related terms and concepts, so that I can improve search accuracy and relevance on the social media platform.'
  'As a product manager, I want to use sentiment analysis to analyze user comments and reviews about our products on social media, so that I can identify areas for improvement and prioritize product development efforts accordingly.'
  'As a customer support manager, I want to implement speech to text for converting voicemail or voice messages from customers into text for easier handling and response management on social media platforms, so that I can improve customer service efficiency.'
  'As a trend analyst, I want to use text categorization to classify social media posts into trending topics or themes (e.g., fashion, technology, politics), so that I can identify popular discussions and trends among users.'
  'As a social network analyst, I want to use unsupervised clustering algorithms to identify and segment communities of users based on their interactions and interests on the platform, so that I can understand community dynamics and foster engagement.'
  'As a content creator, I want to leverage voice recognition capabilities to create and publish audio content (e.g., podcasts, voice notes) directly on social media platforms, so that I can engage with my audience in a more dynamic and personal manner.'
  'As an advertiser, I want to utilize word embedding algorithms to understand the context and meaning behind user-generated content on social networks, so that I can deliver more relevant and targeted advertisements to users.'
)

$clusterName = 'Social Networks'
$domainLabel = 'Social and Urban Studies'
$promptLabel = 'Domain_FSPrompt'

# New fill colour used to highlight the "Social Networks" domain header
# cells (A/B/C), matching xl/styles.xml fgColor rgb="FFE69138" -> BGR int.
$newFillColor = 3707366  # 0x38, 0x91, 0xE6 -> B*65536 + G*256 + R

# ---------------------------------------------------------------------
# 1) Build rows 102-121 by cloning the format of the existing
#    "Social Media" domain block (rows 82-101), then overwriting the
#    values for columns A-F.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 20; $i++) {
    $srcRow = 82 + $i
    $dstRow = 102 + $i

    $ws.Range("A$srcRow`:F$srcRow").Copy()
    $ws.Range("A$dstRow`:F$dstRow").PasteSpecial(-4122)

    $ws.Cells.Item($dstRow, 1).Value = $domainLabel
    $ws.Cells.Item($dstRow, 2).Value = 7.0
    $ws.Cells.Item($dstRow, 3).Value = $clusterName
    $ws.Cells.Item($dstRow, 4).Value = $tasks[$i]
    $ws.Cells.Item($dstRow, 5).Value = $stories[$i]
    $ws.Cells.Item($dstRow, 6).Value = $promptLabel

    # Recolour the domain-header columns (A:C) with the new fill.
    $ws.Range("A$dstRow`:C$dstRow").Interior.Color = $newFillColor
}

# ---------------------------------------------------------------------
# 2) The new block is no longer the tail of the sheet, so every row
#    between the previous last domain row (63) and the new rows
#    (63-121) must show a continuous border down column E, the same
#    way earlier domain boundaries already look in the sheet. Clone
#    the border treatment from an existing pair of rows that already
#    have it (E44 = bottom-only, E45 = top+bottom).
# ---------------------------------------------------------------------
$ws.Range("E44").Copy()
$ws.Range("E63").PasteSpecial(-4122)

$ws.Range("E45").Copy()
$ws.Range("E64:E120").PasteSpecial(-4122)

Write-Output "Added 20 User Story rows for the 'Social Networks' domain."
